# "Generate Report for Handoff"
# Updates the localization-status report: the ebb52c47-... file has been
# handed off (status -> "Ready for handoff") and the c9ca0b1e-... file's
# handoff report also flips to "Ready for handoff"; the latest handback
# timestamps are refreshed and a stale-handback error message is recorded
# for the ebb52c47-... row on the per-language sheets.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$staleHandbackMsg = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54833844ef0c140e87b5919e5bc3891165b32ce8/e2e/ebb52c47-83da-4006-ae86-ec3795f654d8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f418e3a651a16853c8cbd26288d2c114240bb6de/e2e/ebb52c47-83da-4006-ae86-ec3795f654d8.md."

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = "2016-08-22 20:47:46"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $readyForHandoff
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("H3").Value = "2016-08-22 20:47:41"
$wsZhCn.Range("P3").Value = $staleHandbackMsg
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $readyForHandoff
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("H3").Value = "2016-08-22 20:47:46"
$wsDeDe.Range("P3").Value = $staleHandbackMsg
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
